# ASM_Excel.xlsx edit script
# Adds a "Lop:"/"Nhom:" info table at the top of the sheet, adds a third
# group member (SV3) row together with student-id column, and updates a
# couple of task-assignment cells in the work-breakdown table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New "Lop:" / "Nhom:" block in rows 1-2 (replaces the old merged
#    "Nhom: 7" cell that used to live in A2).
# ---------------------------------------------------------------------
$ws.Range("A2").Value2 = ""

$ws.Range("A1").Value2 = "Lớp:"
$ws.Range("B1").Value2 = "UD14311"
$ws.Range("A2").Value2 = "Nhóm:"
$ws.Range("B2").Value2 = 7

$infoLabels = $ws.Range("A1:A2")
$infoLabels.Font.Name = "Times New Roman"
$infoLabels.Font.Size = 13
$infoLabels.Font.Bold = $true
$infoLabels.HorizontalAlignment = -4152   # xlHAlignRight
$infoLabels.VerticalAlignment = -4108     # xlVAlignCenter

$infoValues = $ws.Range("B1:B2")
$infoValues.Font.Name = "Times New Roman"
$infoValues.Font.Size = 13
$infoValues.Font.Bold = $true
$infoValues.HorizontalAlignment = -4131   # xlHAlignLeft
$infoValues.VerticalAlignment = -4108     # xlVAlignCenter

# ---------------------------------------------------------------------
# 2. "Thanh vien nhom" banner (row 9): unmerge A9:C9, move the text to
#    column B and merge B9:D9 is NOT kept (target has no merge there).
# ---------------------------------------------------------------------
$bannerText = $ws.Range("A9").Value2
$ws.Range("A9:C9").UnMerge()
$ws.Range("A9").Value2 = ""
$ws.Range("B9").Value2 = $bannerText
$ws.Range("C9").Value2 = ""
$ws.Range("D9").Value2 = ""

$bannerRange = $ws.Range("B9:D9")
$bannerRange.Font.Name = "Times New Roman"
$bannerRange.Font.Size = 13
$bannerRange.Font.Bold = $false
$bannerRange.HorizontalAlignment = -4131  # xlHAlignLeft (unused -> keep general)
$bannerRange.VerticalAlignment = -4108    # xlVAlignCenter

# Row 9's own alignment in the target keeps only vertical=center (no explicit
# horizontal) - clear horizontal alignment back to general.
$bannerRange.HorizontalAlignment = -4131
$ws.Cells.Item(9, 2).HorizontalAlignment = $null

# ---------------------------------------------------------------------
# 3. SV1..SV4 table (rows 10-13): add a 4th column (D) with student ids
#    and add a new SV3 row. The labels already read SV1:,SV2:,SV3:,SV4:
#    stay in the same rows; just the "PS09108" id and the C-column
#    alignment change.
# ---------------------------------------------------------------------
$ws.Range("D10").Value2 = $ws.Range("D10").Value2   # no-op, keep PS09070
$ws.Range("D11").Value2 = "PS09108"
$ws.Range("D12").Value2 = ""
$ws.Range("D13").Value2 = ""

$memberNames = $ws.Range("C10:C13")
$memberNames.HorizontalAlignment = -4131  # xlHAlignLeft

$memberIds = $ws.Range("D10:D13")
$memberIds.Font.Name = "Times New Roman"
$memberIds.Font.Size = 13
$memberIds.Font.Bold = $false
$memberIds.HorizontalAlignment = -4108    # xlHAlignCenter
$memberIds.VerticalAlignment = -4108      # xlVAlignCenter

# ---------------------------------------------------------------------
# 4. Work breakdown table updates: task "Thiet ke CSDL" (row 18) is now
#    assigned to SV3 alone, and "Bang phan cong cong viec" (row 22) is
#    now assigned to SV4 alone.
# ---------------------------------------------------------------------
$ws.Range("G18").Value2 = "SV3"
$ws.Range("G22").Value2 = "SV4"

# ---------------------------------------------------------------------
# 5. Row heights (rows 9-13 become a fixed 21pt tall; minor re-flow of
#    the rest of the rows after Excel's upgrade/resave).
# ---------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 17.4
$ws.Rows.Item(7).RowHeight = 18.6
$ws.Rows.Item(8).RowHeight = 17.4
$ws.Rows.Item(9).RowHeight = 21
$ws.Rows.Item(10).RowHeight = 21
$ws.Rows.Item(11).RowHeight = 21
$ws.Rows.Item(12).RowHeight = 21
$ws.Rows.Item(13).RowHeight = 21
for ($r = 15; $r -le 22; $r++) {
  $ws.Rows.Item($r).RowHeight = 23.4
}
$ws.Rows.Item(23).RowHeight = 23.4
$ws.Rows.Item(24).RowHeight = 16.2
$ws.Rows.Item(25).RowHeight = 16.2
$ws.Rows.Item(26).RowHeight = 16.2

# ---------------------------------------------------------------------
# 6. Column widths.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 9.21875
$ws.Columns.Item(2).ColumnWidth = 30.6640625
$ws.Columns.Item(3).ColumnWidth = 26.109375
$ws.Columns.Item(4).ColumnWidth = 20.33203125
$ws.Columns.Item(5).ColumnWidth = 24.77734375
$ws.Columns.Item(6).ColumnWidth = 24.33203125
$ws.Columns.Item(7).ColumnWidth = 21.6640625
$ws.Columns.Item(8).ColumnWidth = 20.109375

# ---------------------------------------------------------------------
# 7. View state: scroll so row 13 is at the top and select G21 (matches
#    the author's on-screen state when the file was saved).
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("G21").Select()
